# Fix double-encoded "Â±" (U+00C2 U+00B1) mojibake back to the intended
# plus-minus sign "±" (U+00B1) in every worksheet / every used cell.

$wb = $excel.ActiveWorkbook

$mojibake = [string][char]0x00C2 + [string][char]0x00B1
$correct  = [string][char]0x00B1

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Text
            if ($val -ne $null -and $val.Contains($mojibake)) {
                $cell.Value = $val.Replace($mojibake, $correct)
            }
        }
    }
}
